$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect the latest scrape.
# Rows 8/9 and 38/39 swapped ranking order; other rows only update D (Price) and/or E (Volume 1h).
# D-column price values must stay as text (they are formatted like "64.356.49"), so force
# text number format on cells whose new value would otherwise be auto-parsed as a number.

$ws.Range("D2").Value = '64.356.49'
$ws.Range("E2").Value = '  +4.42%  '
$ws.Range("D3").Value = '2.966.18'
$ws.Range("E3").Value = '  +2.64%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.78'
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.31'
$ws.Range("E6").Value = '  +5.91%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.511'
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = '2.963.34'
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("E10").Value = '  +5.03%  '
$ws.Range("E11").Value = '  +2.08%  '
$ws.Range("E12").Value = '  +3.04%  '
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.31'
$ws.Range("E14").Value = '  +7.30%  '
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '3.458.34'
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").Value = '64.333.64'
$ws.Range("E17").Value = '  +4.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.89'
$ws.Range("E18").Value = '  +4.27%  '
$ws.Range("D19").Value = '2.967.17'
$ws.Range("E19").Value = '  +2.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '444.73'
$ws.Range("E20").Value = '  +2.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.59'
$ws.Range("E21").Value = '  +3.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.676'
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.22'
$ws.Range("E23").Value = '  +4.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.60'
$ws.Range("E24").Value = '  +2.08%  '
$ws.Range("E25").Value = '  +8.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.27'
$ws.Range("E26").Value = '  +2.94%  '
$ws.Range("E27").Value = '  +7.40%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.77'
$ws.Range("E29").Value = '  +9.40%  '
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.20'
$ws.Range("E31").Value = '  +6.81%  '
$ws.Range("E32").Value = '  +1.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.109'
$ws.Range("E33").Value = '  +2.88%  '
$ws.Range("E34").Value = '  +3.43%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.974'
$ws.Range("E36").Value = '  +1.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.63'
$ws.Range("E37").Value = '  +3.58%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.07'
$ws.Range("E38").Value = '  +5.47%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.10'
$ws.Range("E39").Value = '  +7.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.01'
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '43.98'
$ws.Range("E41").Value = '  +13.51%  '
$ws.Range("E42").Value = '  +2.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.293'
$ws.Range("E43").Value = '  +8.74%  '
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '387.37'
$ws.Range("E45").Value = '  +14.72%  '
$ws.Range("D46").Value = '2.768.60'
$ws.Range("E46").Value = '  +3.07%  '
$ws.Range("E47").Value = '  +4.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.90'
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000221'
$ws.Range("E50").Value = '  +12.24%  '
$ws.Range("E51").Value = '  +2.67%  '
